$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the Price column (D) stores values as plain text even though many
# look numeric (Excel would otherwise auto-convert "214.81" style strings
# to real numbers). Values that parse as a normal decimal number get a
# leading apostrophe so COM keeps them as text, matching the workbook's
# original inline-string storage. Values already containing two '.'
# separators (e.g. "29.939.63") never parse as numbers, so no apostrophe
# is needed there.

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "29.939.63"
$ws.Range("E2").Value = "  +0.80%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.633.87"
$ws.Range("E3").Value = "  +1.77%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.25%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "'214.81"
$ws.Range("E5").Value = "  +1.17%  "

# --- Row 6 (XRP) ---
$ws.Range("E6").Value = "  +0.00%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.26%  "

# --- Row 8 (Solana) ---
$ws.Range("D8").Value = "'28.79"
$ws.Range("E8").Value = "  -0.11%  "

# --- Row 9 (Cardano) ---
$ws.Range("E9").Value = "  +0.42%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("E10").Value = "  +0.30%  "

# --- Row 11 (TRON) ---
$ws.Range("D11").Value = "'0.0903"
$ws.Range("E11").Value = "  -0.51%  "

# --- Row 12 (WrappedliquidstakedEther2.0) ---
$ws.Range("D12").Value = "1.867.58"
$ws.Range("E12").Value = "  +1.75%  "

# --- Row 13 (WrappedEther) ---
$ws.Range("D13").Value = "1.633.01"
$ws.Range("E13").Value = "  +2.30%  "

# --- Row 14 (Polygon) ---
$ws.Range("E14").Value = "  +0.21%  "

# --- Row 15 (Chainlink) ---
$ws.Range("D15").Value = "'9.29"
$ws.Range("E15").Value = "  +12.23%  "

# --- Row 16 (WrappedBTC) ---
$ws.Range("D16").Value = "29.959.53"
$ws.Range("E16").Value = "  +0.85%  "

# --- Row 17 (Polkadot) ---
$ws.Range("E17").Value = "  +1.12%  "

# --- Row 18 (Litecoin) ---
$ws.Range("D18").Value = "'64.22"
$ws.Range("E18").Value = "  -0.33%  "

# --- Row 19 (BitcoinCash) ---
$ws.Range("D19").Value = "'241.57"
$ws.Range("E19").Value = "  +0.15%  "

# --- Row 20 (ShibaInu) ---
$ws.Range("E20").Value = "  -0.02%  "

# --- Row 21 (Dai) ---
$ws.Range("E21").Value = "  +0.16%  "

# --- Row 22 (Uniswap) ---
$ws.Range("E22").Value = "  +2.41%  "

# --- Row 23 (Avalanche) ---
$ws.Range("E23").Value = "  +3.30%  "

# --- Row 24 (Toncoin) ---
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  +2.93%  "

# --- Row 25 (Monero) ---
$ws.Range("D25").Value = "'158.08"
$ws.Range("E25").Value = "  +0.74%  "

# --- Row 26 (EthereumClassic) ---
$ws.Range("D26").Value = "'15.48"
$ws.Range("E26").Value = "  -0.38%  "

# --- Row 27 (Stellar) ---
$ws.Range("E27").Value = "  +0.33%  "

# --- Row 28 (Cosmos) ---
$ws.Range("D28").Value = "'6.57"
$ws.Range("E28").Value = "  +1.00%  "

# --- Row 29 (BinanceUSD) ---
$ws.Range("E29").Value = "  +0.18%  "

# --- Row 30 (Hedera) ---
$ws.Range("D30").Value = "'0.0490"
$ws.Range("E30").Value = "  +2.23%  "

# --- Row 31: was PancakeSwap -> now Filecoin ---
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.39"
$ws.Range("E31").Value = "  +4.25%  "

# --- Row 32: was Filecoin -> now PancakeSwap ---
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  +3.34%  "

# --- Row 33 (InternetComputer(DFINITY)) ---
$ws.Range("E33").Value = "  +0.17%  "

# --- Row 34 (Maker) ---
$ws.Range("D34").Value = "1.429.57"
$ws.Range("E34").Value = "  +0.26%  "

# --- Row 35 (LidoDAOToken) ---
$ws.Range("E35").Value = "  +4.88%  "

# --- Row 36 (TrustWalletToken) ---
$ws.Range("E36").Value = "  -1.89%  "

# --- Row 37 (MXToken) ---
$ws.Range("E37").Value = "  -3.05%  "

# --- Row 38 (HuobiToken) ---
$ws.Range("E38").Value = "  +0.07%  "

# --- Row 39 (VeChain) ---
$ws.Range("D39").Value = "'0.0171"
$ws.Range("E39").Value = "  +0.67%  "

# --- Row 40 (Aave) ---
$ws.Range("D40").Value = "'75.45"
$ws.Range("E40").Value = "  +11.18%  "

# --- Row 41 (ImmutableX) ---
$ws.Range("E41").Value = "  -0.13%  "

# --- Row 42 (RenderToken) ---
$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  +0.93%  "

# --- Row 43 (ARBITRUM) ---
$ws.Range("E43").Value = "  +0.65%  "

# --- Row 44 (Kaspa) ---
$ws.Range("D44").Value = "'0.0498"
$ws.Range("E44").Value = "  +0.90%  "

# --- Row 45 (PaxDollar) ---
$ws.Range("E45").Value = "  +0.26%  "

# --- Row 46 (WEMIXToken) ---
$ws.Range("E46").Value = "  +0.66%  "

# --- Row 47 (BitcoinSV) ---
$ws.Range("D47").Value = "'51.23"
$ws.Range("E47").Value = "  -6.25%  "

# --- Row 48 (FraxShare) ---
$ws.Range("D48").Value = "'5.35"
$ws.Range("E48").Value = "  -1.42%  "

# --- Row 49 (RocketPoolETH) ---
$ws.Range("D49").Value = "1.774.60"
$ws.Range("E49").Value = "  +1.75%  "

# --- Row 50 (BabyDogeCoin) ---
$ws.Range("E50").Value = "  +11.05%  "

# --- Row 51 (Quant) ---
$ws.Range("D51").Value = "'90.58"
$ws.Range("E51").Value = "  +4.12%  "
